# DatabaseNormalization.xlsx - typo fix
# - In the "Option B", "Option C", "Option D" and "Option E" sheets, the
#   explanatory paragraph in cell A3 ends the sentence about implicit
#   transitive dependencies with a period even though it introduces a
#   bullet list straight after it; that period should be a colon.
# - Selection/active-cell bookkeeping left behind by the edit (each sheet's
#   active cell moved while the author was reviewing the text).
# - The picture on "Option B" needs its cached bottom-right anchor cell
#   refreshed now that row 3 grew taller to fit the extra wrapped line.

$wb = $excel.ActiveWorkbook

function Fix-Consistency($ws) {
    $cell = $ws.Range("A3")
    $text = $cell.Value()
    $fixed = $text -replace [regex]::Escape("maintain consistency."), "maintain consistency:"
    $cell.Value = $fixed
}

# --- Option A: no text change, just the lingering active-cell selection ---
$wsA = $wb.Worksheets.Item("Option A")
$wsA.Activate()
$wsA.Range("A3").Select()

# --- Option B: typo fix + row height grew (extra wrapped line) + picture ---
$wsB = $wb.Worksheets.Item("Option B")
$wsB.Activate()
Fix-Consistency $wsB
$wsB.Rows.Item(3).RowHeight = 205.65
$wsB.Range("B3").Select()

# Row 3 grew (182.95 -> 205.65pt) to fit the extra wrapped line, which pushes
# the picture's cached bottom-right anchor cell up into row 6 (0-based row 5)
# instead of row 7 (0-based row 6). The picture itself doesn't change size;
# only the cell the bottom-right corner now falls in is recomputed, so the
# row height must be updated *before* nudging the shape to force the anchor
# to resync.
$pic = $wsB.Shapes.Item(1)
$pic.Height = 397.9429

# --- Option C: typo fix ---
$wsC = $wb.Worksheets.Item("Option C")
$wsC.Activate()
Fix-Consistency $wsC
$wsC.Range("B3").Select()

# --- Option D: typo fix ---
$wsD = $wb.Worksheets.Item("Option D")
$wsD.Activate()
Fix-Consistency $wsD
$wsD.Range("B3").Select()

# --- Option E: typo fix ---
$wsE = $wb.Worksheets.Item("Option E")
$wsE.Activate()
Fix-Consistency $wsE
$wsE.Range("B3").Select()

# Restore the originally active sheet (Option E was the selected tab).
$wsE.Activate()
